# "Removed Type from imports"
# The "Commitment Type" column (column G: Sequence, Name, Rule Type, Rule For,
# Generate YTD..., Tag List, Commitment Type, Formula, Entry Type, Rollup,
# Enabled) is removed entirely, shifting Formula/Entry Type/Rollup/Enabled
# one column to the left (H->G, I->H, J->I, K->J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the "Commitment Type" column (matches the author's selection state
# recorded in the saved file: activeCell G1, whole-column sqref) then delete
# it, shifting everything to its right one column to the left.
$ws.Columns("G").Select()
$ws.Columns("G").Delete()
